$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.5253496666666667
$ws.Range("N2").Value = 1.576049
$ws.Range("O2").Value = 0.01837235699974889
$ws.Range("P2").Value = 0.01837235699974889
$ws.Range("Q2").Value = 14.20488007901478
$ws.Range("R2").Value = 127.843920711133
$ws.Range("S2").Value = 0.001303737032974999
$ws.Range("T2").Value = 0.001303737032974999

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.717626000000001
$ws.Range("N3").Value = 8.152878000000001
$ws.Range("O3").Value = 0.09503992908304168
$ws.Range("P3").Value = 0.09503992908304168
$ws.Range("Q3").Value = 73.48163305128068
$ws.Range("R3").Value = 661.334697461526
$ws.Range("S3").Value = 0.006744212250968813
$ws.Range("T3").Value = 0.006744212250968815

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.35159533333334
$ws.Range("N4").Value = 76.05478600000001
$ws.Range("O4").Value = 0.8865877139172095
$ws.Range("P4").Value = 0.8865877139172095
$ws.Range("Q4").Value = 685.4793947175069
$ws.Range("R4").Value = 6169.314552457562
$ws.Range("S4").Value = 0.06291393290639347
$ws.Range("T4").Value = 0.06291393290639348

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.5253496666666667
$ws.Range("N5").Value = 1.576049
$ws.Range("O5").Value = 0.01837235699974889
$ws.Range("P5").Value = 0.01837235699974889
$ws.Range("Q5").Value = 181.5432870887904
$ws.Range("R5").Value = 1633.889583799113
$ws.Range("S5").Value = 0.01666221081410805
$ws.Range("T5").Value = 0.01666221081410805

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.717626000000001
$ws.Range("N6").Value = 8.152878000000001
$ws.Range("O6").Value = 0.09503992908304168
$ws.Range("P6").Value = 0.09503992908304168
$ws.Range("Q6").Value = 939.1207198214541
$ws.Range("R6").Value = 8452.086478393087
$ws.Range("S6").Value = 0.08619336833924809
$ws.Range("T6").Value = 0.08619336833924811

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 25.35159533333334
$ws.Range("N7").Value = 76.05478600000001
$ws.Range("O7").Value = 0.8865877139172095
$ws.Range("P7").Value = 0.8865877139172095
$ws.Range("Q7").Value = 8760.664071532365
$ws.Range("R7").Value = 78845.97664379128
$ws.Range("S7").Value = 0.8040618519816792
$ws.Range("T7").Value = 0.8040618519816793

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.5253496666666667
$ws.Range("N8").Value = 1.576049
$ws.Range("O8").Value = 0.01837235699974889
$ws.Range("P8").Value = 0.01837235699974889
$ws.Range("Q8").Value = 4.428035048953778
$ws.Range("R8").Value = 39.85231544058401
$ws.Range("S8").Value = 0.0004064091526658427
$ws.Range("T8").Value = 0.0004064091526658428

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.717626000000001
$ws.Range("N9").Value = 8.152878000000001
$ws.Range("O9").Value = 0.09503992908304168
$ws.Range("P9").Value = 0.09503992908304168
$ws.Range("Q9").Value = 22.90615934773868
$ws.Range("R9").Value = 206.1554341296481
$ws.Range("S9").Value = 0.002102348492824773
$ws.Range("T9").Value = 0.002102348492824774

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.35159533333334
$ws.Range("N10").Value = 76.05478600000001
$ws.Range("O10").Value = 0.8865877139172095
$ws.Range("P10").Value = 0.8865877139172095
$ws.Range("Q10").Value = 213.6819718477529
$ws.Range("R10").Value = 1923.137746629777
$ws.Range("S10").Value = 0.01961192902913678
$ws.Range("T10").Value = 0.01961192902913679
